$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: paragraphs that contain a grammar-check split ("word" + proofErr +
# "word" runs).  Re-inserting a paragraph's own WordOpenXML collapses the
# run(s) back into a single run (the proofErr markers are transient / not
# preserved across an XML export+import round-trip) while keeping the
# formatting of the paragraph's first run.
# ---------------------------------------------------------------------------
function Collapse-Paragraph([int]$index) {
    $para = $d.Paragraphs.Item($index)
    $rng = $para.Range
    $xml = $rng.WordOpenXML
    $rng.InsertXML($xml)
}

# Find each target paragraph by its (unique) leading text so the script does
# not depend on a hard-coded paragraph index.
function Find-ParagraphIndex([string]$needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.StartsWith($needle)) {
            return $i
        }
    }
    return -1
}

Collapse-Paragraph (Find-ParagraphIndex "Enable Encryption in")
Collapse-Paragraph (Find-ParagraphIndex "Select")
Collapse-Paragraph (Find-ParagraphIndex "Then choose")
Collapse-Paragraph (Find-ParagraphIndex "After encryption")
Collapse-Paragraph (Find-ParagraphIndex "Test the code with sample")
Collapse-Paragraph (Find-ParagraphIndex "Just Save and Test this")

# ---------------------------------------------------------------------------
# Part 2: the hyperlink paragraph.
#   - point it at the new repository URL (this mints a fresh relationship
#     carrying the new target, exactly like Word's own "Edit Hyperlink" UI)
#   - bold the paragraph (both the paragraph-mark run properties and the
#     hyperlink run itself)
# ---------------------------------------------------------------------------
$newUrl = "https://github.com/kohlidevops/aws-lambda-secure-env-variables/blob/main/lambda.py"

$hyperlink = $d.Hyperlinks.Item(1)
$hyperlink.Address = $newUrl
$hyperlink.TextToDisplay = $newUrl

$linkIndex = Find-ParagraphIndex "https://github.com/kohlidevops/aws-lambda-secure-env-variables"
$linkPara = $d.Paragraphs.Item($linkIndex)
$linkRange = $linkPara.Range
$linkXml = $linkRange.WordOpenXML

# Pull out the relationship id Word just assigned to the hyperlink run so the
# hand-edited fragment below keeps pointing at the right (already-updated)
# relationship instead of the package-local numbering used inside the
# exported fragment.
$null = $linkXml -match 'r:id="(rId\d+)"'
$realRelId = $matches[1]

$oldParaRPr = '<w:rFonts w:ascii="Bahnschrift Condensed" w:hAnsi="Bahnschrift Condensed"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:hyperlink'
$newParaRPr = '<w:rFonts w:ascii="Bahnschrift Condensed" w:hAnsi="Bahnschrift Condensed"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:hyperlink'
$linkXml = $linkXml.Replace($oldParaRPr, $newParaRPr)

$oldRunRPr = '<w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Bahnschrift Condensed" w:hAnsi="Bahnschrift Condensed"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr>'
$newRunRPr = '<w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Bahnschrift Condensed" w:hAnsi="Bahnschrift Condensed"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr>'
$linkXml = $linkXml.Replace($oldRunRPr, $newRunRPr)

# make sure the run carries its original rsidR attribute (dropped by the
# WordOpenXML export) and that the hyperlink element references the live
# relationship id rather than the fragment-local one.
$linkXml = $linkXml.Replace('<w:hyperlink r:id="rId5" w:history="1"><w:r w:rsidRPr="005D4F91">', ('<w:hyperlink r:id="' + $realRelId + '" w:history="1"><w:r w:rsidR="009D3740" w:rsidRPr="005D4F91">'))

$linkRange.InsertXML($linkXml)

Write-Host "done"
